$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.730.04"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").Value = "2.519.42"
$ws.Range("E3").Value = "  -2.26%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'304.24"
$ws.Range("E5").Value = "  +0.48%  "

$ws.Range("D6").Value = "'97.72"
$ws.Range("E6").Value = "  +3.20%  "

$ws.Range("D7").Value = "'0.577"
$ws.Range("E7").Value = "  +0.33%  "

$ws.Range("E8").Value = "  +0.17%  "

$ws.Range("D9").Value = "'0.541"
$ws.Range("E9").Value = "  -1.12%  "

$ws.Range("D10").Value = "'36.81"
$ws.Range("E10").Value = "  +1.16%  "

$ws.Range("D11").Value = "'0.0812"
$ws.Range("E11").Value = "  -0.17%  "

$ws.Range("D12").Value = "'7.68"
$ws.Range("E12").Value = "  +0.09%  "

$ws.Range("E13").Value = "  -1.02%  "

$ws.Range("D14").Value = "2.906.01"
$ws.Range("E14").Value = "  -2.13%  "

$ws.Range("D15").Value = "2.522.01"
$ws.Range("E15").Value = "  -1.62%  "

$ws.Range("D16").Value = "'15.10"
$ws.Range("E16").Value = "  +5.80%  "

$ws.Range("D17").Value = "'0.861"
$ws.Range("E17").Value = "  -2.86%  "

$ws.Range("D18").Value = "42.727.36"
$ws.Range("E18").Value = "  -0.71%  "

$ws.Range("D19").Value = "'12.96"
$ws.Range("E19").Value = "  -0.35%  "

$ws.Range("D20").Value = "0.0₃0975"
$ws.Range("E20").Value = "  -2.18%  "

$ws.Range("D21").Value = "'6.44"
$ws.Range("E21").Value = "  -3.50%  "

$ws.Range("D22").Value = "'71.21"
$ws.Range("E22").Value = "  -1.12%  "

$ws.Range("D23").Value = "'251.36"
$ws.Range("E23").Value = "  -1.08%  "

$ws.Range("E24").Value = "  -1.23%  "

$ws.Range("D25").Value = "'2.02"
$ws.Range("E25").Value = "  -5.14%  "

$ws.Range("D26").Value = "'26.89"
$ws.Range("E26").Value = "  -7.05%  "

$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("D28").Value = "'2.33"
$ws.Range("E28").Value = "  +10.82%  "

$ws.Range("D29").Value = "'10.36"
$ws.Range("E29").Value = "  +0.27%  "

$ws.Range("D30").Value = "'38.01"
$ws.Range("E30").Value = "  +1.32%  "

$ws.Range("D31").Value = "'5.98"
$ws.Range("E31").Value = "  -1.05%  "

$ws.Range("D32").Value = "'156.58"
$ws.Range("E32").Value = "  +1.06%  "

$ws.Range("D33").Value = "'0.0792"
$ws.Range("E33").Value = "  -1.48%  "

$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'2.08"
$ws.Range("E34").Value = "  -4.44%  "

$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'3.27"
$ws.Range("E35").Value = "  -4.45%  "

$ws.Range("E36").Value = "  -4.68%  "

$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").Value = "'18.33"
$ws.Range("E37").Value = "  +0.11%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.116"
$ws.Range("E38").Value = "  +1.89%  "

$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "'24.15"
$ws.Range("E39").Value = "  +4.45%  "

$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "'0.119"
$ws.Range("E40").Value = "  -0.92%  "

$ws.Range("E41").Value = "  -6.02%  "

$ws.Range("D42").Value = "'3.40"
$ws.Range("E42").Value = "  -1.31%  "

$ws.Range("D43").Value = "'3.85"
$ws.Range("E43").Value = "  -1.38%  "

$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.18%  "

$ws.Range("D45").Value = "'0.0301"
$ws.Range("E45").Value = "  -3.69%  "

$ws.Range("D46").Value = "2.030.75"
$ws.Range("E46").Value = "  -2.45%  "

$ws.Range("D47").Value = "'85.70"
$ws.Range("E47").Value = "  +0.37%  "

$ws.Range("E48").Value = "  -2.90%  "

$ws.Range("D49").Value = "2.769.45"
$ws.Range("E49").Value = "  -1.85%  "

$ws.Range("D50").Value = "'0.189"
$ws.Range("E50").Value = "  -1.30%  "

$ws.Range("D51").Value = "'101.85"
$ws.Range("E51").Value = "  -4.65%  "
